{"js": "// Replace the 25 \"three-digit x one-digit\" practice equations found in the\n// document's table with their new values. Each table cell holds exactly one\n// equation string (e.g. \"735\u00d72=1470\"); we overwrite them in document order\n// (row by row, left to right) using the old -> new mapping taken from the\n// commit's diff.\n//\n// Positional replacement (by row/column index) is used instead of a global\n// text search-and-replace because a handful of the new equation strings\n// collide with *other* cells' old equation strings (e.g. the new value for\n// the 2nd cell, \"389\u00d77=2723\", is the old value that currently sits in a\n// later cell). A naive \"find this text anywhere and replace the first hit\"\n// approach would therefore clobber the wrong cell once earlier replacements\n// have been written. Addressing cells by (row, column) sidesteps that.\n\nconst oldToNew = [\n  [\"735\u00d72=1470\", \"805\u00d78=6440\"],\n  [\"504\u00d76=3024\", \"389\u00d77=2723\"],\n  [\"292\u00d79=2628\", \"197\u00d78=1576\"],\n  [\"375\u00d77=2625\", \"909\u00d76=5454\"],\n  [\"716\u00d72=1432\", \"163\u00d78=1304\"],\n  [\"559\u00d73=1677\", \"855\u00d73=2565\"],\n  [\"188\u00d74=752\", \"278\u00d74=1112\"],\n  [\"805\u00d76=4830\", \"388\u00d73=1164\"],\n  [\"307\u00d77=2149\", \"879\u00d77=6153\"],\n  [\"749\u00d74=2996\", \"445\u00d73=1335\"],\n  [\"441\u00d76=2646\", \"884\u00d76=5304\"],\n  [\"352\u00d73=1056\", \"573\u00d78=4584\"],\n  [\"605\u00d75=3025\", \"366\u00d73=1098\"],\n  [\"915\u00d74=3660\", \"468\u00d78=3744\"],\n  [\"467\u00d78=3736\", \"897\u00d72=1794\"],\n  [\"258\u00d75=1290\", \"876\u00d76=5256\"],\n  [\"389\u00d77=2723\", \"219\u00d75=1095\"],\n  [\"858\u00d73=2574\", \"899\u00d79=8091\"],\n  [\"385\u00d75=1925\", \"266\u00d75=1330\"],\n  [\"421\u00d75=2105\", \"734\u00d79=6606\"],\n  [\"534\u00d74=2136\", \"881\u00d78=7048\"],\n  [\"595\u00d76=3570\", \"361\u00d79=3249\"],\n  [\"716\u00d74=2864\", \"783\u00d78=6264\"],\n  [\"588\u00d79=5292\", \"363\u00d78=2904\"],\n  [\"955\u00d79=8595\", \"767\u00d74=3068\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table with the practice equations, found none.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's text (one paragraph per cell) so we can find which rows\n// actually contain an equation (several rows in the table are intentionally\n// blank spacer rows).\nconst cellParas = [];\nfor (const row of rows.items) {\n  const rowParas = [];\n  for (const cell of row.cells.items) {\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    rowParas.push(para);\n  }\n  cellParas.push(rowParas);\n}\nawait context.sync();\n\nlet mapIndex = 0;\nfor (const rowParas of cellParas) {\n  const rowHasText = rowParas.some((p) => p.text.trim().length > 0);\n  if (!rowHasText) {\n    continue;\n  }\n  for (const para of rowParas) {\n    if (mapIndex >= oldToNew.length) {\n      break;\n    }\n    const [oldText, newText] = oldToNew[mapIndex];\n    const actual = para.text.trim();\n    if (actual !== oldText) {\n      throw new Error(\n        `Unexpected cell text at position ${mapIndex}: expected \"${oldText}\" but found \"${actual}\"`\n      );\n    }\n    para.getRange().insertText(newText, Word.InsertLocation.replace);\n    mapIndex++;\n  }\n}\nawait context.sync();\n\nif (mapIndex !== oldToNew.length) {\n  throw new Error(`Only replaced ${mapIndex} of ${oldToNew.length} equations.`);\n}\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" practice equations found in the\n# document's table with their new values. Each table cell holds exactly one\n# equation string (e.g. \"735\u00d72=1470\"); we overwrite them in document order\n# (row by row, left to right) using the old -> new mapping taken from the\n# commit's diff.\n#\n# Positional replacement (by row/column index via $table.Cell(r,c)) is used\n# instead of a document-wide Find/Replace because a handful of the new\n# equation strings collide with *other* cells' old equation strings (e.g.\n# the new value for the 2nd cell, \"389\u00d77=2723\", is the old value that\n# currently sits in a later cell). A blind \"find this text anywhere and\n# replace it\" pass would clobber the wrong cell once earlier replacements\n# have already been written. Addressing cells by (row, column) sidesteps\n# that entirely.\n\n$oldToNew = @(\n    ,@(\"735\u00d72=1470\", \"805\u00d78=6440\")\n    ,@(\"504\u00d76=3024\", \"389\u00d77=2723\")\n    ,@(\"292\u00d79=2628\", \"197\u00d78=1576\")\n    ,@(\"375\u00d77=2625\", \"909\u00d76=5454\")\n    ,@(\"716\u00d72=1432\", \"163\u00d78=1304\")\n    ,@(\"559\u00d73=1677\", \"855\u00d73=2565\")\n    ,@(\"188\u00d74=752\", \"278\u00d74=1112\")\n    ,@(\"805\u00d76=4830\", \"388\u00d73=1164\")\n    ,@(\"307\u00d77=2149\", \"879\u00d77=6153\")\n    ,@(\"749\u00d74=2996\", \"445\u00d73=1335\")\n    ,@(\"441\u00d76=2646\", \"884\u00d76=5304\")\n    ,@(\"352\u00d73=1056\", \"573\u00d78=4584\")\n    ,@(\"605\u00d75=3025\", \"366\u00d73=1098\")\n    ,@(\"915\u00d74=3660\", \"468\u00d78=3744\")\n    ,@(\"467\u00d78=3736\", \"897\u00d72=1794\")\n    ,@(\"258\u00d75=1290\", \"876\u00d76=5256\")\n    ,@(\"389\u00d77=2723\", \"219\u00d75=1095\")\n    ,@(\"858\u00d73=2574\", \"899\u00d79=8091\")\n    ,@(\"385\u00d75=1925\", \"266\u00d75=1330\")\n    ,@(\"421\u00d75=2105\", \"734\u00d79=6606\")\n    ,@(\"534\u00d74=2136\", \"881\u00d78=7048\")\n    ,@(\"595\u00d76=3570\", \"361\u00d79=3249\")\n    ,@(\"716\u00d74=2864\", \"783\u00d78=6264\")\n    ,@(\"588\u00d79=5292\", \"363\u00d78=2904\")\n    ,@(\"955\u00d79=8595\", \"767\u00d74=3068\")\n)\n\n$d = $word.ActiveDocument\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected a table with the practice equations, found none.\"\n}\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$mapIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    # Collect this row's cell text (trim the trailing cell-mark chars: CR + BEL).\n    $rowTexts = @()\n    for ($c = 1; $c -le $colCount; $c++) {\n        $raw = $table.Cell($r, $c).Range.Text\n        $rowTexts += $raw.TrimEnd([char]13, [char]7)\n    }\n\n    $rowHasText = $false\n    foreach ($cellText in $rowTexts) {\n        if ($cellText.Trim().Length -gt 0) {\n            $rowHasText = $true\n        }\n    }\n    if (-not $rowHasText) {\n        continue\n    }\n\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($mapIndex -ge $oldToNew.Length) {\n            break\n        }\n        $pair = $oldToNew[$mapIndex]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $actual = $rowTexts[$c - 1].Trim()\n        if ($actual -ne $oldText) {\n            throw \"Unexpected cell text at row $r col $c (position $mapIndex): expected [$oldText] but found [$actual]\"\n        }\n        $table.Cell($r, $c).Range.Text = $newText\n        $mapIndex++\n    }\n}\n\nif ($mapIndex -ne $oldToNew.Length) {\n    throw \"Only replaced $mapIndex of $($oldToNew.Length) equations.\"\n}\n\nWrite-Output \"Replaced $mapIndex equations.\"\n"}
